$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of COVID overview data to append (dates must stay plain text,
# matching the existing sheet's convention, so force text format first).
$newRows = @(
    @("2021-06-14", "overview", "K02000001", "United Kingdom", 4573419, 7742, 3, 127907),
    @("2021-06-15", "overview", "K02000001", "United Kingdom", 4581006, 7673, 10, 127917),
    @("2021-06-16", "overview", "K02000001", "United Kingdom", 4589814, 9055, 9, 127926),
    @("2021-06-17", "overview", "K02000001", "United Kingdom", 4600623, 11007, 19, 127945),
    @("2021-06-18", "overview", "K02000001", "United Kingdom", 4610893, 10476, 11, 127956)
)

$startRow = 307
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]

    # Columns A-D are text (date stored as plain text, like the rest of the sheet)
    $textRange = $ws.Range("A" + $r + ":D" + $r)
    $textRange.NumberFormat = "@"
    $ws.Cells.Item($r, 1).Value = $rowData[0]
    $ws.Cells.Item($r, 2).Value = $rowData[1]
    $ws.Cells.Item($r, 3).Value = $rowData[2]
    $ws.Cells.Item($r, 4).Value = $rowData[3]

    # Columns E-H are numbers
    $ws.Cells.Item($r, 5).Value = $rowData[4]
    $ws.Cells.Item($r, 6).Value = $rowData[5]
    $ws.Cells.Item($r, 7).Value = $rowData[6]
    $ws.Cells.Item($r, 8).Value = $rowData[7]
}
